$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 data rows (old "MuSCs-as-sender" rows); new dataset only has
# FAPs and MuSCs as sending clusters (ECs dropped as sender), so the table shrinks
# from 10 rows (9 data rows) to 7 rows (6 data rows).
$ws.Rows("8:10").Delete() | Out-Null

# Refresh all data rows (2-7) with the recomputed TPM-based NATMI values.

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.39906333333333
$ws.Range("H2").Value = 52.19719000000001
$ws.Range("I2").Value = 0.9351306508759385
$ws.Range("J2").Value = 0.9351306508759385
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 39.93062416217668
$ws.Range("R2").Value = 359.37561745959
$ws.Range("S2").Value = 0.01480085534149662
$ws.Range("T2").Value = 0.01480085534149662

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.39906333333333
$ws.Range("H3").Value = 52.19719000000001
$ws.Range("I3").Value = 0.9351306508759385
$ws.Range("J3").Value = 0.9351306508759385
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 111.5917106666667
$ws.Range("N3").Value = 334.775132
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 1941.591241364342
$ws.Range("R3").Value = 17474.32117227908
$ws.Range("S3").Value = 0.7196784848399915
$ws.Range("T3").Value = 0.7196784848399918

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.39906333333333
$ws.Range("H4").Value = 52.19719000000001
$ws.Range("I4").Value = 0.9351306508759385
$ws.Range("J4").Value = 0.9351306508759385
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 541.3289901243012
$ws.Range("R4").Value = 4871.960911118711
$ws.Range("S4").Value = 0.2006513106944502
$ws.Range("T4").Value = 0.2006513106944502

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.206960666666667
$ws.Range("H5").Value = 3.620882
$ws.Range("I5").Value = 0.06486934912406146
$ws.Range("J5").Value = 0.06486934912406146
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 2.769959039511333
$ws.Range("R5").Value = 24.929631355602
$ws.Range("S5").Value = 0.001026724823513085
$ws.Range("T5").Value = 0.001026724823513085

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.206960666666667
$ws.Range("H6").Value = 3.620882
$ws.Range("I6").Value = 0.06486934912406146
$ws.Range("J6").Value = 0.06486934912406146
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 111.5917106666667
$ws.Range("N6").Value = 334.775132
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("Q6").Value = 134.6868055007137
$ws.Range("R6").Value = 1212.181249506424
$ws.Range("S6").Value = 0.04992358537968036
$ws.Range("T6").Value = 0.04992358537968038

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.206960666666667
$ws.Range("H7").Value = 3.620882
$ws.Range("I7").Value = 0.06486934912406146
$ws.Range("J7").Value = 0.06486934912406146
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("Q7").Value = 37.55160759457089
$ws.Range("R7").Value = 337.964468351138
$ws.Range("S7").Value = 0.013919038920868
$ws.Range("T7").Value = 0.01391903892086801
